$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from row 255 (column A date style) down to new rows 256-269
$ws.Range("A255").Copy($ws.Range("A256:A269"))

$ws.Range("A256").Value = 44330
$ws.Range("B256").Value = 9
$ws.Range("C256").Value = 66
$ws.Range("D256").Value = 91.19927040583676

$ws.Range("A257").Value = 44331
$ws.Range("B257").Value = 5
$ws.Range("C257").Value = 59
$ws.Range("D257").Value = 81.52662051430862

$ws.Range("A258").Value = 44332
$ws.Range("B258").Value = 13
$ws.Range("C258").Value = 61
$ws.Range("D258").Value = 84.29023476903093

$ws.Range("A259").Value = 44333
$ws.Range("B259").Value = 5
$ws.Range("C259").Value = 56
$ws.Range("D259").Value = 77.38119913222513

$ws.Range("A260").Value = 44334
$ws.Range("B260").Value = 4
$ws.Range("C260").Value = 50
$ws.Range("D260").Value = 69.09035636805815

$ws.Range("A261").Value = 44335
$ws.Range("B261").Value = 0
$ws.Range("C261").Value = 50
$ws.Range("D261").Value = 69.09035636805815

$ws.Range("A262").Value = 44336
$ws.Range("B262").Value = 18
$ws.Range("C262").Value = 54
$ws.Range("D262").Value = 74.61758487750279

$ws.Range("A263").Value = 44337
$ws.Range("B263").Value = 7
$ws.Range("C263").Value = 52
$ws.Range("D263").Value = 71.85397062278047

$ws.Range("A264").Value = 44338
$ws.Range("B264").Value = 6
$ws.Range("C264").Value = 53
$ws.Range("D264").Value = 73.23577775014164

$ws.Range("A265").Value = 44339
$ws.Range("B265").Value = 0
$ws.Range("C265").Value = 40
$ws.Range("D265").Value = 55.27228509444652

$ws.Range("A266").Value = 44340
$ws.Range("B266").Value = 7
$ws.Range("C266").Value = 42
$ws.Range("D266").Value = 58.03589934916884

$ws.Range("A267").Value = 44341
$ws.Range("B267").Value = 3
$ws.Range("C267").Value = 41
$ws.Range("D267").Value = 56.65409222180768

$ws.Range("A268").Value = 44342
$ws.Range("B268").Value = 2
$ws.Range("C268").Value = 43
$ws.Range("D268").Value = 59.41770647653001

$ws.Range("A269").Value = 44343
$ws.Range("B269").Value = 1
$ws.Range("C269").Value = 26
$ws.Range("D269").Value = 35.92698531139023
